# Insert a new data row before row 37 (pushes existing rows 37-64 down to 38-65)
# and populate it with the new weekly record, matching the author's commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(37).Insert()

$ws.Cells.Item(37, 1).Value = 6
$ws.Cells.Item(37, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(37, 3).Value = "Metropolitana"
$ws.Cells.Item(37, 4).Value = 45049
$ws.Cells.Item(37, 5).Value = 13
$ws.Cells.Item(37, 6).Value = 100112035
$ws.Cells.Item(37, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 18000
$ws.Cells.Item(37, 12).Value = 20000
$ws.Cells.Item(37, 13).Value = 18850
$ws.Cells.Item(37, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(37, 16).Value = 1257
$ws.Cells.Item(37, 17).Value = 15
$ws.Cells.Item(37, 18).Value = "Hortaliza"
